$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.65427347373088
$ws.Range("C2").Value = 8.344836826775365
$ws.Range("D2").Value = 12.012965196457
$ws.Range("F2").Value = 26.97652056924061
$ws.Range("G2").Value = 23.29348078567289
$ws.Range("H2").Value = 12.79410509837156
$ws.Range("I2").Value = 17.39621189606612
$ws.Range("J2").Value = 11.2541287738502
$ws.Range("N2").Value = 16.01916465888374
$ws.Range("O2").Value = 18.75932125869748

$ws.Range("B3").Value = 12.05982313196583
$ws.Range("C3").Value = 7.846818874029967
$ws.Range("D3").Value = 11.93121086638853
$ws.Range("F3").Value = 26.94930551006954
$ws.Range("G3").Value = 23.23616245126769
$ws.Range("H3").Value = 12.83295452102592
$ws.Range("I3").Value = 17.49644792272704
$ws.Range("J3").Value = 11.23320048125566
$ws.Range("N3").Value = 16.03961755552672
$ws.Range("O3").Value = 18.80309558906082

$ws.Range("B4").Value = 11.67978812193205
$ws.Range("C4").Value = 7.522956823737785
$ws.Range("D4").Value = 11.88308798599628
$ws.Range("F4").Value = 26.94042887965885
$ws.Range("G4").Value = 23.21049007013245
$ws.Range("H4").Value = 12.85919712544585
$ws.Range("I4").Value = 17.56198020200425
$ws.Range("J4").Value = 11.22262592083332
$ws.Range("N4").Value = 16.05384375254147
$ws.Range("O4").Value = 18.83490592589779

$ws.Range("B5").Value = 11.52132449386295
$ws.Range("C5").Value = 7.386466985174074
$ws.Range("D5").Value = 11.8640162211613
$ws.Range("F5").Value = 26.93878405492825
$ws.Range("G5").Value = 23.20242815961331
$ws.Range("H5").Value = 12.87049115142452
$ws.Range("I5").Value = 17.58968681410518
$ws.Range("J5").Value = 11.21889192026213
$ws.Range("N5").Value = 16.06006128543374
$ws.Range("O5").Value = 18.84910549350425

$ws.Range("B6").Value = 11.49480052143988
$ws.Range("C6").Value = 7.363531483694705
$ws.Range("D6").Value = 11.86088238521692
$ws.Range("F6").Value = 26.93863010861232
$ws.Range("G6").Value = 23.20123451998855
$ws.Range("H6").Value = 12.87240272872975
$ws.Range("I6").Value = 17.59434795032244
$ws.Range("J6").Value = 11.21830671320788
$ws.Range("N6").Value = 16.06111910677415
$ws.Range("O6").Value = 18.85153789109537

$ws.Range("B7").Value = 11.67766532060809
$ws.Range("C7").Value = 7.521134299237757
$ws.Range("D7").Value = 11.882828574507
$ws.Range("F7").Value = 26.94039870832771
$ws.Range("G7").Value = 23.21037162302777
$ws.Range("H7").Value = 12.85934701248016
$ws.Range("I7").Value = 17.56234980843773
$ws.Range("J7").Value = 11.22257323011225
$ws.Range("N7").Value = 16.05392590177537
$ws.Range("O7").Value = 18.83509242456364

$ws.Range("B8").Value = 12.45253126214907
$ws.Range("C8").Value = 8.176888587629479
$ws.Range("D8").Value = 11.98435619450281
$ws.Range("F8").Value = 26.96551310248116
$ws.Range("G8").Value = 23.27174571611152
$ws.Range("H8").Value = 12.80700411669117
$ws.Range("I8").Value = 17.42994524863426
$ws.Range("J8").Value = 11.24644237463988
$ws.Range("N8").Value = 16.02587110066818
$ws.Range("O8").Value = 18.77338869432563

$ws.Range("B9").Value = 13.84591383624992
$ws.Range("C9").Value = 9.318410964498693
$ws.Range("D9").Value = 12.19912531249973
$ws.Range("F9").Value = 27.07672144279875
$ws.Range("G9").Value = 23.46723115444847
$ws.Range("H9").Value = 12.72334768997851
$ws.Range("I9").Value = 17.20197487269365
$ws.Range("J9").Value = 11.31115179937256
$ws.Range("N9").Value = 15.98405588970303
$ws.Range("O9").Value = 18.69168106159965

$ws.Range("B10").Value = 14.78530271093547
$ws.Range("C10").Value = 10.06822720362179
$ws.Range("D10").Value = 12.36536657961192
$ws.Range("F10").Value = 27.19581929248173
$ws.Range("G10").Value = 23.65585332623706
$ws.Range("H10").Value = 12.67350452445191
$ws.Range("I10").Value = 17.05385061053497
$ws.Range("J10").Value = 11.36937376489043
$ws.Range("N10").Value = 15.9613347613916
$ws.Range("O10").Value = 18.65579350227019

$ws.Range("B11").Value = 15.1931737343152
$ws.Range("C11").Value = 10.3899726267745
$ws.Range("D11").Value = 12.44258125977984
$ws.Range("F11").Value = 27.25800054214475
$ws.Range("G11").Value = 23.75118849635573
$ws.Range("H11").Value = 12.65336129681963
$ws.Range("I11").Value = 16.99068359648017
$ws.Range("J11").Value = 11.39812006380566
$ws.Range("N11").Value = 15.95272535065847
$ws.Range("O11").Value = 18.64474181368699

$ws.Range("B12").Value = 15.34474923295887
$ws.Range("C12").Value = 10.50902502723185
$ws.Range("D12").Value = 12.47202768388509
$ws.Range("F12").Value = 27.28268480310896
$ws.Range("G12").Value = 23.78863411693295
$ws.Range("H12").Value = 12.64609815679345
$ws.Range("I12").Value = 16.96737145448267
$ws.Range("J12").Value = 11.40932487075061
$ws.Range("N12").Value = 15.94971256892199
$ws.Range("O12").Value = 18.64131706508707

$ws.Range("B13").Value = 15.31223376138817
$ws.Range("C13").Value = 10.48350887297176
$ws.Range("D13").Value = 12.46567701842497
$ws.Range("F13").Value = 27.27731824209706
$ws.Range("G13").Value = 23.78051019598363
$ws.Range("H13").Value = 12.64764617392336
$ws.Range("I13").Value = 16.97236508046019
$ws.Range("J13").Value = 11.40689761679346
$ws.Range("N13").Value = 15.95035043544814
$ws.Range("O13").Value = 18.64202080632506

$ws.Range("B14").Value = 15.20570198667854
$ws.Range("C14").Value = 10.39982297493372
$ws.Range("D14").Value = 12.44499981434524
$ws.Range("F14").Value = 27.26000861480079
$ws.Range("G14").Value = 23.75424239008152
$ws.Range("H14").Value = 12.65275644136408
$ws.Range("I14").Value = 16.98875350081638
$ws.Range("J14").Value = 11.3990355279533
$ws.Range("N14").Value = 15.95247253372755
$ws.Range("O14").Value = 18.64444481091265

$ws.Range("B15").Value = 15.14007152949251
$ws.Range("C15").Value = 10.3482001337722
$ws.Range("D15").Value = 12.43236072114965
$ws.Range("F15").Value = 27.24955368970942
$ws.Range("G15").Value = 23.73832683297978
$ws.Range("H15").Value = 12.65593413926012
$ws.Range("I15").Value = 16.99887109018781
$ws.Range("J15").Value = 11.39426115927149
$ws.Range("N15").Value = 15.95380457506106
$ws.Range("O15").Value = 18.64602864465021

$ws.Range("B16").Value = 14.75824878879097
$ws.Range("C16").Value = 10.04681068113161
$ws.Range("D16").Value = 12.3603505677308
$ws.Range("F16").Value = 27.19191555100832
$ws.Range("G16").Value = 23.64981255495909
$ws.Range("H16").Value = 12.67487192099282
$ws.Range("I16").Value = 17.05806368907112
$ws.Range("J16").Value = 11.36754012628791
$ws.Range("N16").Value = 15.96193207614801
$ws.Range("O16").Value = 18.6566220527477

$ws.Range("B17").Value = 14.51896841600656
$ws.Range("C17").Value = 9.856958884299807
$ws.Range("D17").Value = 12.3165666503839
$ws.Range("F17").Value = 27.15859733033287
$ws.Range("G17").Value = 23.59793519615445
$ws.Range("H17").Value = 12.68713835510837
$ws.Range("I17").Value = 17.09545726448189
$ws.Range("J17").Value = 11.35172263587186
$ws.Range("N17").Value = 15.96735963979547
$ws.Range("O17").Value = 18.66447301663779

$ws.Range("B18").Value = 14.37951253712603
$ws.Range("C18").Value = 9.745939837906818
$ws.Range("D18").Value = 12.29153435323008
$ws.Range("F18").Value = 27.14018776399596
$ws.Range("G18").Value = 23.56899549427239
$ws.Range("H18").Value = 12.69443184057979
$ws.Range("I18").Value = 17.11736173442349
$ws.Range("J18").Value = 11.34283801564721
$ws.Range("N18").Value = 15.97064399420842
$ws.Range("O18").Value = 18.66948498668993

$ws.Range("B19").Value = 14.33198379057749
$ws.Range("C19").Value = 9.708038028688753
$ws.Range("D19").Value = 12.28308545142623
$ws.Range("F19").Value = 27.13408450901221
$ws.Range("G19").Value = 23.55935211923921
$ws.Range("H19").Value = 12.69694216943049
$ws.Range("I19").Value = 17.1248462916171
$ws.Range("J19").Value = 11.33986662258724
$ws.Range("N19").Value = 15.97178396642271
$ws.Range("O19").Value = 18.6712671278589

$ws.Range("B20").Value = 14.54463003321805
$ws.Range("C20").Value = 9.877357463307522
$ws.Range("D20").Value = 12.32121204682275
$ws.Range("F20").Value = 27.16206614597421
$ws.Range("G20").Value = 23.60336478812437
$ws.Range("H20").Value = 12.68580791828234
$ws.Range("I20").Value = 17.09143558846838
$ws.Range("J20").Value = 11.35338441253506
$ws.Range("N20").Value = 15.96676504783317
$ws.Range("O20").Value = 18.66358588713994

$ws.Range("B21").Value = 15.23707156886024
$ws.Range("C21").Value = 10.42447916314777
$ws.Range("D21").Value = 12.45106777067589
$ws.Range("F21").Value = 27.26506211306184
$ws.Range("G21").Value = 23.76192162257011
$ws.Range("H21").Value = 12.65124552929772
$ws.Range("I21").Value = 16.98392331440235
$ws.Range("J21").Value = 11.40133619998852
$ws.Range("N21").Value = 15.9518425142087
$ws.Range("O21").Value = 18.64371217501303

$ws.Range("B22").Value = 15.67282936113643
$ws.Range("C22").Value = 10.76581862121595
$ws.Range("D22").Value = 12.53713165517757
$ws.Range("F22").Value = 27.33899917079995
$ws.Range("G22").Value = 23.87336874594208
$ws.Range("H22").Value = 12.63078286373295
$ws.Range("I22").Value = 16.91720158506829
$ws.Range("J22").Value = 11.43453301148844
$ws.Range("N22").Value = 15.94353148609637
$ws.Range("O22").Value = 18.63515564566504

$ws.Range("B23").Value = 15.44181604570739
$ws.Range("C23").Value = 10.58512526115548
$ws.Range("D23").Value = 12.49109562198607
$ws.Range("F23").Value = 27.29893638243632
$ws.Range("G23").Value = 23.813181133937
$ws.Range("H23").Value = 12.64150942011057
$ws.Range("I23").Value = 16.95248740476921
$ws.Range("J23").Value = 11.41664734273469
$ws.Range("N23").Value = 15.9478356116362
$ws.Range("O23").Value = 18.63931637426667

$ws.Range("B24").Value = 14.53303429821034
$ws.Range("C24").Value = 9.868141092388386
$ws.Range("D24").Value = 12.3191114264057
$ws.Range("F24").Value = 27.16049557117332
$ws.Range("G24").Value = 23.60090731004642
$ws.Range("H24").Value = 12.68640865702922
$ws.Range("I24").Value = 17.09325252237494
$ws.Range("J24").Value = 11.35263247164952
$ws.Range("N24").Value = 15.96703335201635
$ws.Range("O24").Value = 18.66398540620601

$ws.Range("B25").Value = 13.48332008857521
$ws.Range("C25").Value = 9.025165657012321
$ws.Range("D25").Value = 12.1394580517626
$ws.Range("F25").Value = 27.04003527973728
$ws.Range("G25").Value = 23.4063726758324
$ws.Range("H25").Value = 12.74394153584219
$ws.Range("I25").Value = 17.26025043395008
$ws.Range("J25").Value = 11.29175187894579
$ws.Range("N25").Value = 15.99395958656353
$ws.Range("O25").Value = 18.70955742920438
